# Update the cryptos.xlsx "Sat Jan 7 17:39:17 UTC 2023" symbol-list refresh.
# All touched cells (D = Price, E = Volume(1h)) are stored as plain text
# (inlineStr) in the workbook, so we force Text formatting before writing
# the new values to keep Excel from reinterpreting them as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "D2"  = "261.08"
    "D3"  = "27.09"
    "E3"  = "0.85%"
    "D4"  = "4.705"
    "D5"  = "0.06181"
    "E5"  = "3.35%"
    "D6"  = "6.684"
    "E6"  = "0.61%"
    "D7"  = "0.8522"
    "E7"  = "-0.51%"
    "D8"  = "0.9128"
    "E8"  = "-0.97%"
    "E9"  = "1.39%"
    "D10" = "0.04672"
    "E10" = "2.75%"
    "D11" = "0.07088"
    "E11" = "1.18%"
    "D12" = "0.03146"
    "E12" = "3.12%"
    "D13" = "0.09039"
    "D14" = "0.001527"
    "E14" = "-0.38%"
    "D15" = "0.0006163"
    "E15" = "1.66%"
    "D16" = "0.006128"
    "E16" = "-0.92%"
    "D17" = "3.457"
    "E17" = "0.35%"
    "D18" = "3.177"
    "E18" = "0.84%"
    "D19" = "2.153"
    "E19" = "-0.03%"
    "E20" = "-0.88%"
    "D22" = "4.078"
    "E22" = "1.29%"
    "D23" = "0.04213"
    "E23" = "-0.33%"
    "E24" = "-0.01%"
    "E25" = "-5.60%"
    "E26" = "0.07%"
    "D40" = "0.03888"
    "E40" = "1.48%"
    "D41" = "0.1111"
    "E41" = "-0.10%"
    "E42" = "9.04%"
    "E43" = "8.53%"
    "E44" = "-9.71%"
    "D45" = "0.00005160"
    "E45" = "0.97%"
    "E46" = "0.07%"
    "D48" = "0.1677"
    "E48" = "6.37%"
    "E49" = "0.07%"
    "E50" = "0.07%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Keep the cell formatted as Text so the string value round-trips
    # byte-for-byte (e.g. "0.07%" instead of being parsed into a numeric
    # percentage, "261.08" instead of a float losing its text nature).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
